$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# This document (OC0202 angivPassword.docx) gets a handful of changes:
#  1. The title paragraph "OC0202- angivPassword" + each of its runs gets an
#     explicit English (US) language mark (<w:lang w:val="en-US"/>).
#  2. The bordered (underline) paragraph right below the title also gets the
#     same language mark on its paragraph mark.
#  3. The two following blank paragraphs (previously totally empty <w:p/>)
#     get a <w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr> too.
#  4. "UC02 Delete User" is retranslated to Danish: "UC02 Slet Bruger",
#     spell-check flagging "Slet" and "Bruger" as the (foreign-looking)
#     words, same as "angivPassword" elsewhere in the doc.
#  5. "Klienten k" is renamed to "Patienten p" (merged back into one run).
# --------------------------------------------------------------------------

$langRPr = '<w:rPr><w:lang w:val="en-US"/></w:rPr>'

function Assert-ParaText($para, $expected, $label) {
    $actual = $para.Range.Text
    if ($actual -ne $expected) {
        Write-Host "WARNING: paragraph '$label' text mismatch. Expected:[$expected] Actual:[$actual]"
    }
}

# 1. Title paragraph: "OC0202- angivPassword"
$p1 = $d.Paragraphs.Item(1)
Assert-ParaText $p1 "OC0202- angivPassword`r" "title"
$p1xml = '<w:p w14:paraId="43331344" w14:textId="166988F6" w:rsidR="00334F33" w:rsidRPr="0051282E" w:rsidRDefault="0051282E" w:rsidP="001259E4">' + `
  '<w:pPr><w:pStyle w:val="Titel"/>' + $langRPr + '</w:pPr>' + `
  '<w:r w:rsidRPr="0051282E">' + $langRPr + '<w:t>OC0</w:t></w:r>' + `
  '<w:r w:rsidR="001B0100">' + $langRPr + '<w:t>2</w:t></w:r>' + `
  '<w:r w:rsidRPr="0051282E">' + $langRPr + '<w:t>0</w:t></w:r>' + `
  '<w:r w:rsidR="003E5E79">' + $langRPr + '<w:t>2</w:t></w:r>' + `
  '<w:r w:rsidRPr="0051282E">' + $langRPr + '<w:t xml:space="preserve">- </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r w:rsidR="003E5E79">' + $langRPr + '<w:t>angivPassword</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>'
[void]$p1.Range.InsertXML($p1xml)

# 2. Bordered blank paragraph right after the title
$p2 = $d.Paragraphs.Item(2)
Assert-ParaText $p2 "`r" "border paragraph"
$p2xml = '<w:p w14:paraId="581FC920" w14:textId="77777777" w:rsidR="001259E4" w:rsidRPr="0051282E" w:rsidRDefault="001259E4" w:rsidP="001259E4">' + `
  '<w:pPr><w:pBdr><w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/></w:pBdr>' + $langRPr + '</w:pPr>' + `
  '</w:p>'
[void]$p2.Range.InsertXML($p2xml)

# 3 & 4. The two empty paragraphs right after that
$p3 = $d.Paragraphs.Item(3)
Assert-ParaText $p3 "`r" "empty paragraph 1"
$p3xml = '<w:p w14:paraId="305B3794" w14:textId="366FFFD5" w:rsidR="001259E4" w:rsidRPr="0051282E" w:rsidRDefault="001259E4" w:rsidP="001259E4">' + `
  '<w:pPr>' + $langRPr + '</w:pPr>' + `
  '</w:p>'
[void]$p3.Range.InsertXML($p3xml)

$p4 = $d.Paragraphs.Item(4)
$p4xml = '<w:p w14:paraId="01B66AC5" w14:textId="77777777" w:rsidR="001259E4" w:rsidRPr="0051282E" w:rsidRDefault="001259E4" w:rsidP="001259E4">' + `
  '<w:pPr>' + $langRPr + '</w:pPr>' + `
  '</w:p>'
[void]$p4.Range.InsertXML($p4xml)

# 8. "UC02 Delete User" -> "UC02 Slet Bruger"
$p8 = $d.Paragraphs.Item(8)
$p8xml = '<w:p w14:paraId="07289A78" w14:textId="76F00434" w:rsidR="001259E4" w:rsidRPr="001B0100" w:rsidRDefault="004D7F30">' + `
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r w:rsidRPr="001B0100"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>UC0</w:t></w:r>' + `
  '<w:r w:rsidR="00F52C31" w:rsidRPr="001B0100"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">2 </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r w:rsidR="00F52C31" w:rsidRPr="001B0100"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Slet</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r w:rsidR="00F52C31" w:rsidRPr="001B0100"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r w:rsidR="00F52C31" w:rsidRPr="001B0100"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Bruger</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>'
[void]$p8.Range.InsertXML($p8xml)

# 11. "Klienten" + " " + "k" -> single run "Patienten p"
$p11 = $d.Paragraphs.Item(11)
$p11xml = '<w:p w14:paraId="67A4AB2A" w14:textId="4AAA63F9" w:rsidR="00075D1B" w:rsidRDefault="00075D1B" w:rsidP="004D7F30">' + `
  '<w:r><w:t>Patienten p</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> er logget ind på brugeren</w:t></w:r>' + `
  '</w:p>'
[void]$p11.Range.InsertXML($p11xml)

Write-Host "Done applying edits"
